# Updates crypto price/volume figures (and swaps OKB/ONDO row order) to match
# the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.301.46"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "'3.572.66"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'599.93"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "'134.83"
$ws.Range("E6").Value = "  -3.77%  "

$ws.Range("D7").Value = "'3.571.37"
$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("E10").Value = "  -1.54%  "

$ws.Range("D11").Value = "'7.19"
$ws.Range("E11").Value = "  +1.96%  "

$ws.Range("D12").Value = "'0.390"
$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("D13").Value = "'4.185.66"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").Value = "'0.0000184"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").Value = "'27.36"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").Value = "'3.575.97"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").Value = "'64.510.46"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("D19").Value = "'9.99"
$ws.Range("E19").Value = "  -2.81%  "

$ws.Range("D20").Value = "'14.43"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").Value = "'5.86"
$ws.Range("E21").Value = "  -0.43%  "

$ws.Range("D22").Value = "'392.09"
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("D23").Value = "'0.580"
$ws.Range("E23").Value = "  +1.13%  "

$ws.Range("D24").Value = "'3.717.88"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").Value = "'74.15"
$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").Value = "'0.0000114"
$ws.Range("E27").Value = "  -2.17%  "

$ws.Range("D28").Value = "'7.90"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").Value = "'1.62"
$ws.Range("E29").Value = "  +28.28%  "

$ws.Range("D30").Value = "'8.68"
$ws.Range("E30").Value = "  +4.86%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("D33").Value = "'3.577.00"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "'24.26"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("D37").Value = "'171.01"
$ws.Range("E37").Value = "  +1.56%  "

$ws.Range("D38").Value = "'6.97"
$ws.Range("E38").Value = "  -1.42%  "

$ws.Range("D39").Value = "'5.13"
$ws.Range("E39").Value = "  +2.45%  "

$ws.Range("D40").Value = "'1.55"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("D41").Value = "'0.0820"
$ws.Range("E41").Value = "  +1.73%  "

$ws.Range("D42").Value = "'0.829"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").Value = "'26.36"
$ws.Range("E43").Value = "  -1.50%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'43.14"
$ws.Range("E44").Value = "  +0.53%  "

$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "'1.25"
$ws.Range("E45").Value = "  +5.08%  "

$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").Value = "'4.47"
$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("D48").Value = "'1.67"
$ws.Range("E48").Value = "  -1.42%  "

$ws.Range("D49").Value = "'6.95"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("D50").Value = "'2.443.26"
$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("D51").Value = "'0.0268"
$ws.Range("E51").Value = "  +1.31%  "
